$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Restricciones_del_follower
# Regenerated "follower" restriction rows (rows 2-6, columns A:F).
# Columns B, D, E, F contain numbers stored as *text* in the source workbook,
# so force Text formatting before assigning them (otherwise Excel would
# auto-convert the numeric-looking strings into real numbers).
# ---------------------------------------------------------------------------
# NOTE: worksheet names are resolved by index, not name. This workbook has
# two sheets whose names differ only by case ("Vector_bf" / "Vector_BF"),
# and name lookups here resolve case-insensitively (always landing on the
# first match) -- so we pin every sheet to its 1-based tab position.
$wsFollower = $wb.Worksheets.Item(3)

# Force text storage on every cell that holds a numeric-looking value.
# (Two single-area calls instead of one comma-joined multi-area range: the
# multi-area form only applies the format to its first area.)
$wsFollower.Range("B2:B6").NumberFormat = "@"
$wsFollower.Range("D2:F6").NumberFormat = "@"

$wsFollower.Range("A2").Value = "8.600000000000001 - y_1"
$wsFollower.Range("B2").Value = "-8.600000000000001"
$wsFollower.Range("C2").Value = "J_0_L0_v"
$wsFollower.Range("D2").Value = "0.75"
$wsFollower.Range("E2").Value = "8.5"
$wsFollower.Range("F2").Value = "0"

$wsFollower.Range("A3").Value = "-8.600000000000001 + y_1"
$wsFollower.Range("B3").Value = "4.600000000000001"
$wsFollower.Range("C3").Value = "J_0_L0_v"
$wsFollower.Range("D3").Value = "0.19"
$wsFollower.Range("E3").Value = "0"
$wsFollower.Range("F3").Value = "6.2"

$wsFollower.Range("A4").Value = "-5.000000000000002 - 2x + y_1 + 4y_2"
$wsFollower.Range("B4").Value = "-10.999999999999998"
$wsFollower.Range("C4").Value = "J_0_LP_v"
$wsFollower.Range("D4").Value = "0.24"
$wsFollower.Range("E4").Value = "0"
$wsFollower.Range("F4").Value = "3.0"

$wsFollower.Range("A5").Value = "-65.41 + 8x + y_1"
$wsFollower.Range("B5").Value = "17.4"
$wsFollower.Range("C5").Value = "J_Ne_L0_v"
$wsFollower.Range("D5").Value = "0.92"
$wsFollower.Range("E5").Value = "0"
$wsFollower.Range("F5").Value = "6.3"

$wsFollower.Range("A6").Value = "-7.400000000000002 - 2x - 2y_1"
$wsFollower.Range("B6").Value = "-19.400000000000002"
$wsFollower.Range("C6").Value = "J_Ne_L0_v"
$wsFollower.Range("D6").Value = "0.82"
$wsFollower.Range("E6").Value = "0"
$wsFollower.Range("F6").Value = "0.3"

# ---------------------------------------------------------------------------
# Sheet: Punto_modificado (x, y_1, y_2)
# ---------------------------------------------------------------------------
$wsPunto = $wb.Worksheets.Item(4)
$wsPunto.Range("A2:C2").NumberFormat = "@"
$wsPunto.Range("A2").Value = "7.1"
$wsPunto.Range("B2").Value = "8.600000000000001"
$wsPunto.Range("C2").Value = "2.65"

# ---------------------------------------------------------------------------
# Sheet: Vector_bf
# ---------------------------------------------------------------------------
$wsBf = $wb.Worksheets.Item(5)
$wsBf.Range("A2:A3").NumberFormat = "@"
$wsBf.Range("A2").Value = "2.04"
$wsBf.Range("A3").Value = "-0.96"

# ---------------------------------------------------------------------------
# Sheet: Vector_BF
# ---------------------------------------------------------------------------
$wsBF = $wb.Worksheets.Item(6)
$wsBF.Range("A2:A4").NumberFormat = "@"
$wsBF.Range("A2").Value = "1.0"
$wsBF.Range("A3").Value = "11.5"
$wsBF.Range("A4").Value = "-2.0"
